# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (interest count) figures in column F for the
# "展览" and "全部类型" sheets to reflect the latest scrape.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 41
$wsExhibition.Range("F4").Value = 232
$wsExhibition.Range("F5").Value = 2758
$wsExhibition.Range("F6").Value = 1930
$wsExhibition.Range("F7").Value = 374
$wsExhibition.Range("F9").Value = 980
$wsExhibition.Range("F11").Value = 21

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 41
$wsAll.Range("F4").Value = 232
$wsAll.Range("F5").Value = 2758
$wsAll.Range("F6").Value = 1930
$wsAll.Range("F7").Value = 374
$wsAll.Range("F10").Value = 980
$wsAll.Range("F12").Value = 21
